# The exam taker filled in their answers (column C) for worksheet "3", which are
# compared against the correct answers already stored in column B by the existing
# D-column formulas. Most rows get the same value as column B (correct answer),
# a handful of rows are left blank (unanswered), and a handful get a different,
# wrong answer (plus a yellow highlight on a few of the answered cells, matching
# the highlighting already used elsewhere in the workbook for flagged answers).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3")

$ws.Range("C2").Value = "Use the 802.1x protocol"
$ws.Range("C3").Value = "Place a front-end web server in a demilitarized zone that only handles external web traffic"
$ws.Range("C5").Value = "Source IP address"
$ws.Range("C7").Value = "Test password strength, brute-force encrypted or hashed passwords, and crack passwords via dictionary attacks"
$ws.Range("C8").Value = "Grey-box"
$ws.Range("C9").Value = "BBProxy"
$ws.Range("C10").Value = "Salting"
$ws.Range("C14").Value = "Default Credential"
$ws.Range("C15").Value = "Nothing because the password file does not contain the passwords themselves"
$ws.Range("C16").Value = "2"
$ws.Range("C17").Value = "Watering Hole"
$ws.Range("C18").Value = "Rules of Engagement"
$ws.Range("C19").Value = "IP Spoofing"
$ws.Range("C20").Value = "Display of the contents of the passwd file"
$ws.Range("C22").Value = "C"
$ws.Range("C24").Value = "Open-source intelligence"
$ws.Range("C25").Value = "DNS spoofing"
$ws.Range("C26").Value = "Trap"
$ws.Range("C27").Value = "Nslookup"
$ws.Range("C28").Value = "Sniffers operate on Layer 2 of the OSI model"
$ws.Range("C31").Value = "Windows"
$ws.Range("C32").Value = "Confidentiality"
$ws.Range("C33").Value = "Open mail relay"
$ws.Range("C34").Value = "Port scanning, banner grabbing, service identification"
$ws.Range("C36").Value = "http-methods"
$ws.Range("C37").Value = "Fuzz Testing"
$ws.Range("C38").Value = "SYN"
$ws.Range("C39").Value = "Perform a vulnerability scan of the system"
$ws.Range("C41").Value = "Brute force"
$ws.Range("C42").Value = "Nikto"
$ws.Range("C43").Value = "Physical layer"
$ws.Range("C46").Value = "Markov Chain"
$ws.Range("C49").Value = "A web server facing the Internet, an application server on the internal network, a database server on the internal network"
$ws.Range("C49").Interior.Color = 65535
$ws.Range("C50").Value = "RSA is asymmetric, which is used to create a public/private key pair; AES is symmetric, which is used to encrypt data"
$ws.Range("C50").Interior.Color = 65535
$ws.Range("C51").Value = "Nessus"
$ws.Range("C52").Value = "T5"
$ws.Range("C56").Value = "NTP"
$ws.Range("C57").Value = "tcp.port ==21"
$ws.Range("C58").Value = "Encrypt the backup tapes and transport them in a lockbox"
$ws.Range("C60").Value = "The tester only partially knows the internal structure"
$ws.Range("C61").Value = "The data collection speeds, data processing speed, or enrolment time"
$ws.Range("C62").Value = "He is scanning from 192.168.1.64 to 192.168.1.78 because of the mask /28 and the servers are not in that range"
$ws.Range("C63").Value = "An inability to access any website"
$ws.Range("C64").Value = "Social Engineering"
$ws.Range("C65").Value = "Chosen-plaintext"
$ws.Range("C66").Value = "Not informing the employees that they are going to be monitored could be an invasion of privacy"
$ws.Range("C67").Value = "Business Impact Analysis (BIA)"
$ws.Range("C68").Value = "Gray Hat"
$ws.Range("C69").Value = "Clickjacking"
$ws.Range("C70").Value = "Use the 3-2-1 backup rule"
$ws.Range("C71").Value = "ESP transport mode"
$ws.Range("C72").Value = "Demilitarized Zone"
$ws.Range("C73").Value = "The host is likely a printer"
$ws.Range("C74").Value = "Steganography"
$ws.Range("C75").Value = "Wired Equivalent Privacy (WEP)"
$ws.Range("C76").Value = "Double quotation"
$ws.Range("C78").Value = "John should write to jacksmith@gmail.com to verify the identity of Jack"
$ws.Range("C78").Interior.Color = 65535
$ws.Range("C79").Value = "Bollards"
$ws.Range("C80").Value = "Remote access policy"
$ws.Range("C81").Value = "DNSSEC"
$ws.Range("C83").Value = "If (source matches 10.10.10.0/24 and destination matches 10.20.20.1 and port matches 443) then permit"
$ws.Range("C84").Value = "Confidentiality, integrity, and availability"
$ws.Range("C85").Value = "Automatic and manual testing should be used together to better cover potential problems"
$ws.Range("C86").Value = "Public Key"
$ws.Range("C87").Value = "Use hping"
$ws.Range("C88").Value = "Creating a botnet"
$ws.Range("C92").Value = "Use a scan tool like Nessus"
$ws.Range("C94").Value = "Corporate Espionage"
$ws.Range("C96").Value = "This is the process of sniffing through the switch"
$ws.Range("C96").Interior.Color = 65535
$ws.Range("C97").Value = "Kernel mode rootkits"
$ws.Range("C98").Value = "Unicode characters"
$ws.Range("C99").Value = "Use Tor network with multi-node"
$ws.Range("C100").Value = "Use password salting"
$ws.Range("C101").Value = "It sends a request packet to all the network elements, asking for the MAC address from a specific IP"
$ws.Range("C102").Value = "Banner grabbing"
$ws.Range("C103").Value = "Height/Weight"
$ws.Range("C104").Value = "Dictionary-attack"
$ws.Range("C105").Value = "SSL/TLS uses both asymmetric and symmetric encryption"
$ws.Range("C106").Value = "CHNTPW"
$ws.Range("C107").Value = "Ettercap"
$ws.Range("C108").Value = "XSS"
$ws.Range("C109").Value = "Reconnaissance"
$ws.Range("C110").Value = "Evil Twin"
$ws.Range("C111").Value = "Attempts by attackers to access passwords stored on the employee's computer"
$ws.Range("C111").Interior.Color = 65535
$ws.Range("C112").Value = "Threat"
$ws.Range("C114").Value = "Rainbow Table Attack"
$ws.Range("C115").Value = "Reconnaissance"
$ws.Range("C116").Value = "TCP SYN (Stealth) Scan"
$ws.Range("C116").Interior.Color = 65535
$ws.Range("C117").Value = "Macro Virus"
$ws.Range("C118").Value = "Confidentiality"
$ws.Range("C119").Value = "He will create a SPAN entry on the spoofed root bridge and redirect traffic to his computer"
$ws.Range("C120").Value = "False Positives and False Negatives"
$ws.Range("C121").Value = "Code Emulation"
$ws.Range("C122").Value = "Single sign-on"
$ws.Range("C124").Value = "Zero-Day"
$ws.Range("C125").Value = "Macro virus"

# Make worksheet "3" the active sheet/tab and restore the selection the author
# left behind, matching the saved view state.
$ws.Activate()
$ws.Range("C2:C126").Select()
